$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ash_data")

# Fill in rows 64-73 (A=sequential index continuing from 60, C=0, D=1, E=0, F=3)
for ($row = 64; $row -le 73; $row++) {
    $idx = $row - 3
    $ws.Cells.Item($row, 1).Value = $idx
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 1
    $ws.Cells.Item($row, 5).Value = 0
    $ws.Cells.Item($row, 6).Value = 3
}

# Update the view: scroll position and selection
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("A52:F73").Select() | Out-Null
